$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update passenger name split: "Neeha Rathna" / "Janjanam" -> "Neeha" / "Rathna Janjanam"
$ws.Range("F5").Value = "Neeha"
$ws.Range("G5").Value = "Rathna Janjanam"

# Update the active selection on the sheet
$ws.Range("F9").Select()
